# Apply the "modern response figures, starting insurance and invest" edit.
#
# Summary of changes:
#  - survey sheet: row 4 ("insurance_choice" question) used the wrong
#    select type "select_one invest_options"; fix it to
#    "select_one insurance_options". Widen column A to fit the new text
#    and move the active selection there.
#  - choices sheet: the insurance answer rows (4-5) used the wrong
#    list_name "insurance_choice"; fix it to "insurance_options" so it
#    matches the corrected survey sheet. Make "choices" the active tab
#    (instead of "settings") with the selection on A6.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")
$settings = $wb.Worksheets.Item("settings")

# Fix the survey sheet's insurance question type.
$survey.Range("A4").Value = "select_one insurance_options"

# Fix the choices sheet's insurance list_name entries.
$choices.Range("A4").Value = "insurance_options"
$choices.Range("A5").Value = "insurance_options"

# Widen column A on the survey sheet to fit the longer text.
$survey.Columns.Item(1).ColumnWidth = 25.17

# Update the survey sheet's own (non-active) selection to A4.
$survey.Range("A4").Select() | Out-Null

# Make "choices" the active sheet/tab, with its selection on A6.
[void]$choices.Activate()
$choices.Range("A6").Select() | Out-Null
